$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.092.55"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.533.12"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'606.74"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'143.81"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("D7").Value = "3.531.63"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("D11").Value = "'8.06"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").Value = "4.128.57"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("E14").Value = "  -4.65%  "
$ws.Range("D15").Value = "'30.39"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").Value = "3.529.43"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "66.230.28"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "'15.03"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").Value = "'425.56"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").Value = "'78.96"
$ws.Range("D25").Value = "3.676.61"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").Value = "'8.01"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").Value = "'2.48"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'1.49"
$ws.Range("E32").Value = "  -5.83%  "
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "'25.35"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "3.521.33"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'5.60"
$ws.Range("E39").Value = "  -6.06%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'170.64"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'0.0860"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.892"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'5.17"
$ws.Range("E44").Value = "  -5.37%  "
$ws.Range("E45").Value = "  -10.20%  "
$ws.Range("D46").Value = "'45.30"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("E47").Value = "  -9.84%  "
$ws.Range("D48").Value = "'25.86"
$ws.Range("E48").Value = "  -8.40%  "
$ws.Range("D49").Value = "'2.40"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -4.06%  "
